# Add album/track stats to the Artists sheet, matching the wireframe update.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Artists")
$ws3 = $wb.Worksheets.Item("Tracks")

# New header cells for the stats columns.
$ws1.Range("B1").Value = "Albums"
$ws1.Range("C1").Value = "Tracks"

# Match the existing bold / left-aligned header style (A1) across B1:C1.
$ws1.Range("A1:C1").HorizontalAlignment = -4131
$ws1.Range("A1:C1").Font.Bold = $true

# Left-align the rest of the table (A2:C5), matching the existing look.
$ws1.Range("A2:C5").HorizontalAlignment = -4131

# Album / track counts per artist.
$ws1.Range("B2").Value = 1
$ws1.Range("C2").Value = 9

$ws1.Range("B3").Value = 1
$ws1.Range("C3").Value = 13

$ws1.Range("B4").Value = 1
$ws1.Range("C4").Value = 5

$ws1.Range("B5").Value = 2
$ws1.Range("C5").Value = 29

# Tracks sheet keeps its existing selection, but is no longer the active tab.
[void]$ws3.Range("A2").Select()

# Artists becomes the active sheet/tab, with the selection left on C6.
[void]$ws1.Range("C6").Select()
